# Actualización automática 2025-09-26 12:55:09
#
# A new client, "BARROS YUNGA DIEGO VINICIO", is inserted alphabetically as
# the second client row (row 3) on both worksheets, pushing all the
# following client rows down by one. The new client's figures are all 0.
# The summary/totals row at the bottom shifts down with it, and on the
# "VENTAS POR GRUPO" sheet the "X de N" counter text updates its
# denominator from 11 to 12 clients.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Push rows 3..13 down to 4..14, leaving a blank row 3 for the new client.
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B3").Value = "BARROS YUNGA DIEGO VINICIO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(3, $col).Value = 0
}

# The totals row (now row 14) counts "X de 11" -> "X de 12" clients.
$ws1.Range("C14").Value = "0 de 12"
$ws1.Range("D14").Value = "0 de 12"
$ws1.Range("E14").Value = "0 de 12"
$ws1.Range("F14").Value = "0 de 12"
$ws1.Range("G14").Value = "0 de 12"
$ws1.Range("H14").Value = "0 de 12"
$ws1.Range("I14").Value = "0 de 12"
$ws1.Range("J14").Value = "0 de 12"
$ws1.Range("K14").Value = "0 de 12"
$ws1.Range("L14").Value = "0 de 12"
$ws1.Range("M14").Value = "3 de 12"
$ws1.Range("N14").Value = "0 de 12"
$ws1.Range("O14").Value = "0 de 12"
$ws1.Range("P14").Value = "1 de 12"
$ws1.Range("Q14").Value = "0 de 12"
$ws1.Range("R14").Value = "0 de 12"

# ---------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Push rows 3..13 down to 4..14, leaving a blank row 3 for the new client.
$ws2.Rows.Item(3).Insert()

$ws2.Range("A3").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B3").Value = "BARROS YUNGA DIEGO VINICIO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(3, $col).Value = 0
}
